$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Cells.Clear()

# --- Column B labels: establishes shared-string order for indices 0-27 ---
$ws.Range("B2").Value = "HKL"
$ws.Range("B3").Value = "Spiral5"
$ws.Range("B4").Value = "RotRing OmegaMax-90"
$ws.Range("B5").Value = "Equal Angle"
$ws.Range("B6").Value = "Tilt Rotate"
$ws.Range("B7").Value = "CLR"
$ws.Range("B8").Value = "Rizzie Hex"
$ws.Range("B9").Value = "Thomas Hex"
$ws.Range("B10").Value = "Tilt Rotate_Partial"
$ws.Range("B11").Value = "RotRing OmegaMax-60"
$ws.Range("B12").Value = "Equal Angle_Partial"
$ws.Range("B13").Value = "Rizzie Hex_Partial"
$ws.Range("B14").Value = "ND Single"
$ws.Range("B15").Value = "RD Single"
$ws.Range("B16").Value = "TD Single"
$ws.Range("B17").Value = "Morris Single"
$ws.Range("B18").Value = "Ring Perpendicular to ND"
$ws.Range("B19").Value = "Ring Perpendicular to RD"
$ws.Range("B20").Value = "Ring Perpendicular to TD"
$ws.Range("B21").Value = "OffsetFTD"
$ws.Range("B22").Value = "OffsetATD"
$ws.Range("B23").Value = "OffsetF45"
$ws.Range("B24").Value = "OffsetA45"
$ws.Range("B25").Value = "OffsetFRD"
$ws.Range("B26").Value = "OffsetARD"
$ws.Range("B27").Value = "Gaussian Quadrature"
$ws.Range("B28").Value = "Michael-CCHex"
$ws.Range("B29").Value = "Michael-SNHex"

# --- Row 2 HKL/pair labels: establishes shared-string order for indices 28-48 ---
$ws.Range("C2").Value = "[4, 2, 0]"
$ws.Range("D2").Value = "[4, 0, 0]"
$ws.Range("E2").Value = "[2, 0, 0]"
$ws.Range("F2").Value = "[2, 2, 0]"
$ws.Range("G2").Value = "[3, 3, 3]"
$ws.Range("H2").Value = "[1, 1, 1]"
$ws.Range("I2").Value = "[2, 2, 2]"
$ws.Range("J2").Value = "[3, 3, 1]"
$ws.Range("K2").Value = "[3, 1, 1]"
$ws.Range("L2").Value = "[4, 2, 2]"
$ws.Range("M2").Value = "[5, 1, 1]"
$ws.Range("N2").Value = "1Pair-A"
$ws.Range("O2").Value = "1Pair-B"
$ws.Range("P2").Value = "2Pairs-A"
$ws.Range("Q2").Value = "2Pairs-B"
$ws.Range("R2").Value = "3Pairs-A"
$ws.Range("S2").Value = "3Pairs-B"
$ws.Range("T2").Value = "3Pairs-C"
$ws.Range("U2").Value = "4Pairs"
$ws.Range("V2").Value = "5A4F"
$ws.Range("W2").Value = "MaxUnique"

# --- Row 1 numeric header (B1:W1 = 0..21) ---
$ws.Range("B1").Value = 0
$ws.Range("C1").Value = 1
$ws.Range("D1").Value = 2
$ws.Range("E1").Value = 3
$ws.Range("F1").Value = 4
$ws.Range("G1").Value = 5
$ws.Range("H1").Value = 6
$ws.Range("I1").Value = 7
$ws.Range("J1").Value = 8
$ws.Range("K1").Value = 9
$ws.Range("L1").Value = 10
$ws.Range("M1").Value = 11
$ws.Range("N1").Value = 12
$ws.Range("O1").Value = 13
$ws.Range("P1").Value = 14
$ws.Range("Q1").Value = 15
$ws.Range("R1").Value = 16
$ws.Range("S1").Value = 17
$ws.Range("T1").Value = 18
$ws.Range("U1").Value = 19
$ws.Range("V1").Value = 20
$ws.Range("W1").Value = 21

# --- Column A numeric index (A2:A29 = 0..27) ---
$ws.Range("A2").Value = 0
$ws.Range("A3").Value = 1
$ws.Range("A4").Value = 2
$ws.Range("A5").Value = 3
$ws.Range("A6").Value = 4
$ws.Range("A7").Value = 5
$ws.Range("A8").Value = 6
$ws.Range("A9").Value = 7
$ws.Range("A10").Value = 8
$ws.Range("A11").Value = 9
$ws.Range("A12").Value = 10
$ws.Range("A13").Value = 11
$ws.Range("A14").Value = 12
$ws.Range("A15").Value = 13
$ws.Range("A16").Value = 14
$ws.Range("A17").Value = 15
$ws.Range("A18").Value = 16
$ws.Range("A19").Value = 17
$ws.Range("A20").Value = 18
$ws.Range("A21").Value = 19
$ws.Range("A22").Value = 20
$ws.Range("A23").Value = 21
$ws.Range("A24").Value = 22
$ws.Range("A25").Value = 23
$ws.Range("A26").Value = 24
$ws.Range("A27").Value = 25
$ws.Range("A28").Value = 26
$ws.Range("A29").Value = 27

# --- Data grid (C3:W29 = 1) ---
$ws.Range("C3").Value = 1
$ws.Range("D3").Value = 1
$ws.Range("E3").Value = 1
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 1
$ws.Range("H3").Value = 1
$ws.Range("I3").Value = 1
$ws.Range("J3").Value = 1
$ws.Range("K3").Value = 1
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 1
$ws.Range("N3").Value = 1
$ws.Range("O3").Value = 1
$ws.Range("P3").Value = 1
$ws.Range("Q3").Value = 1
$ws.Range("R3").Value = 1
$ws.Range("S3").Value = 1
$ws.Range("T3").Value = 1
$ws.Range("U3").Value = 1
$ws.Range("V3").Value = 1
$ws.Range("W3").Value = 1
$ws.Range("C4").Value = 1
$ws.Range("D4").Value = 1
$ws.Range("E4").Value = 1
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 1
$ws.Range("H4").Value = 1
$ws.Range("I4").Value = 1
$ws.Range("J4").Value = 1
$ws.Range("K4").Value = 1
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 1
$ws.Range("N4").Value = 1
$ws.Range("O4").Value = 1
$ws.Range("P4").Value = 1
$ws.Range("Q4").Value = 1
$ws.Range("R4").Value = 1
$ws.Range("S4").Value = 1
$ws.Range("T4").Value = 1
$ws.Range("U4").Value = 1
$ws.Range("V4").Value = 1
$ws.Range("W4").Value = 1
$ws.Range("C5").Value = 1
$ws.Range("D5").Value = 1
$ws.Range("E5").Value = 1
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 1
$ws.Range("H5").Value = 1
$ws.Range("I5").Value = 1
$ws.Range("J5").Value = 1
$ws.Range("K5").Value = 1
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 1
$ws.Range("N5").Value = 1
$ws.Range("O5").Value = 1
$ws.Range("P5").Value = 1
$ws.Range("Q5").Value = 1
$ws.Range("R5").Value = 1
$ws.Range("S5").Value = 1
$ws.Range("T5").Value = 1
$ws.Range("U5").Value = 1
$ws.Range("V5").Value = 1
$ws.Range("W5").Value = 1
$ws.Range("C6").Value = 1
$ws.Range("D6").Value = 1
$ws.Range("E6").Value = 1
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 1
$ws.Range("H6").Value = 1
$ws.Range("I6").Value = 1
$ws.Range("J6").Value = 1
$ws.Range("K6").Value = 1
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 1
$ws.Range("N6").Value = 1
$ws.Range("O6").Value = 1
$ws.Range("P6").Value = 1
$ws.Range("Q6").Value = 1
$ws.Range("R6").Value = 1
$ws.Range("S6").Value = 1
$ws.Range("T6").Value = 1
$ws.Range("U6").Value = 1
$ws.Range("V6").Value = 1
$ws.Range("W6").Value = 1
$ws.Range("C7").Value = 1
$ws.Range("D7").Value = 1
$ws.Range("E7").Value = 1
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 1
$ws.Range("H7").Value = 1
$ws.Range("I7").Value = 1
$ws.Range("J7").Value = 1
$ws.Range("K7").Value = 1
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 1
$ws.Range("N7").Value = 1
$ws.Range("O7").Value = 1
$ws.Range("P7").Value = 1
$ws.Range("Q7").Value = 1
$ws.Range("R7").Value = 1
$ws.Range("S7").Value = 1
$ws.Range("T7").Value = 1
$ws.Range("U7").Value = 1
$ws.Range("V7").Value = 1
$ws.Range("W7").Value = 1
$ws.Range("C8").Value = 1
$ws.Range("D8").Value = 1
$ws.Range("E8").Value = 1
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 1
$ws.Range("H8").Value = 1
$ws.Range("I8").Value = 1
$ws.Range("J8").Value = 1
$ws.Range("K8").Value = 1
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 1
$ws.Range("N8").Value = 1
$ws.Range("O8").Value = 1
$ws.Range("P8").Value = 1
$ws.Range("Q8").Value = 1
$ws.Range("R8").Value = 1
$ws.Range("S8").Value = 1
$ws.Range("T8").Value = 1
$ws.Range("U8").Value = 1
$ws.Range("V8").Value = 1
$ws.Range("W8").Value = 1
$ws.Range("C9").Value = 1
$ws.Range("D9").Value = 1
$ws.Range("E9").Value = 1
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 1
$ws.Range("H9").Value = 1
$ws.Range("I9").Value = 1
$ws.Range("J9").Value = 1
$ws.Range("K9").Value = 1
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 1
$ws.Range("N9").Value = 1
$ws.Range("O9").Value = 1
$ws.Range("P9").Value = 1
$ws.Range("Q9").Value = 1
$ws.Range("R9").Value = 1
$ws.Range("S9").Value = 1
$ws.Range("T9").Value = 1
$ws.Range("U9").Value = 1
$ws.Range("V9").Value = 1
$ws.Range("W9").Value = 1
$ws.Range("C10").Value = 1
$ws.Range("D10").Value = 1
$ws.Range("E10").Value = 1
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 1
$ws.Range("H10").Value = 1
$ws.Range("I10").Value = 1
$ws.Range("J10").Value = 1
$ws.Range("K10").Value = 1
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 1
$ws.Range("N10").Value = 1
$ws.Range("O10").Value = 1
$ws.Range("P10").Value = 1
$ws.Range("Q10").Value = 1
$ws.Range("R10").Value = 1
$ws.Range("S10").Value = 1
$ws.Range("T10").Value = 1
$ws.Range("U10").Value = 1
$ws.Range("V10").Value = 1
$ws.Range("W10").Value = 1
$ws.Range("C11").Value = 1
$ws.Range("D11").Value = 1
$ws.Range("E11").Value = 1
$ws.Range("F11").Value = 1
$ws.Range("G11").Value = 1
$ws.Range("H11").Value = 1
$ws.Range("I11").Value = 1
$ws.Range("J11").Value = 1
$ws.Range("K11").Value = 1
$ws.Range("L11").Value = 1
$ws.Range("M11").Value = 1
$ws.Range("N11").Value = 1
$ws.Range("O11").Value = 1
$ws.Range("P11").Value = 1
$ws.Range("Q11").Value = 1
$ws.Range("R11").Value = 1
$ws.Range("S11").Value = 1
$ws.Range("T11").Value = 1
$ws.Range("U11").Value = 1
$ws.Range("V11").Value = 1
$ws.Range("W11").Value = 1
$ws.Range("C12").Value = 1
$ws.Range("D12").Value = 1
$ws.Range("E12").Value = 1
$ws.Range("F12").Value = 1
$ws.Range("G12").Value = 1
$ws.Range("H12").Value = 1
$ws.Range("I12").Value = 1
$ws.Range("J12").Value = 1
$ws.Range("K12").Value = 1
$ws.Range("L12").Value = 1
$ws.Range("M12").Value = 1
$ws.Range("N12").Value = 1
$ws.Range("O12").Value = 1
$ws.Range("P12").Value = 1
$ws.Range("Q12").Value = 1
$ws.Range("R12").Value = 1
$ws.Range("S12").Value = 1
$ws.Range("T12").Value = 1
$ws.Range("U12").Value = 1
$ws.Range("V12").Value = 1
$ws.Range("W12").Value = 1
$ws.Range("C13").Value = 1
$ws.Range("D13").Value = 1
$ws.Range("E13").Value = 1
$ws.Range("F13").Value = 1
$ws.Range("G13").Value = 1
$ws.Range("H13").Value = 1
$ws.Range("I13").Value = 1
$ws.Range("J13").Value = 1
$ws.Range("K13").Value = 1
$ws.Range("L13").Value = 1
$ws.Range("M13").Value = 1
$ws.Range("N13").Value = 1
$ws.Range("O13").Value = 1
$ws.Range("P13").Value = 1
$ws.Range("Q13").Value = 1
$ws.Range("R13").Value = 1
$ws.Range("S13").Value = 1
$ws.Range("T13").Value = 1
$ws.Range("U13").Value = 1
$ws.Range("V13").Value = 1
$ws.Range("W13").Value = 1
$ws.Range("C14").Value = 1
$ws.Range("D14").Value = 1
$ws.Range("E14").Value = 1
$ws.Range("F14").Value = 1
$ws.Range("G14").Value = 1
$ws.Range("H14").Value = 1
$ws.Range("I14").Value = 1
$ws.Range("J14").Value = 1
$ws.Range("K14").Value = 1
$ws.Range("L14").Value = 1
$ws.Range("M14").Value = 1
$ws.Range("N14").Value = 1
$ws.Range("O14").Value = 1
$ws.Range("P14").Value = 1
$ws.Range("Q14").Value = 1
$ws.Range("R14").Value = 1
$ws.Range("S14").Value = 1
$ws.Range("T14").Value = 1
$ws.Range("U14").Value = 1
$ws.Range("V14").Value = 1
$ws.Range("W14").Value = 1
$ws.Range("C15").Value = 1
$ws.Range("D15").Value = 1
$ws.Range("E15").Value = 1
$ws.Range("F15").Value = 1
$ws.Range("G15").Value = 1
$ws.Range("H15").Value = 1
$ws.Range("I15").Value = 1
$ws.Range("J15").Value = 1
$ws.Range("K15").Value = 1
$ws.Range("L15").Value = 1
$ws.Range("M15").Value = 1
$ws.Range("N15").Value = 1
$ws.Range("O15").Value = 1
$ws.Range("P15").Value = 1
$ws.Range("Q15").Value = 1
$ws.Range("R15").Value = 1
$ws.Range("S15").Value = 1
$ws.Range("T15").Value = 1
$ws.Range("U15").Value = 1
$ws.Range("V15").Value = 1
$ws.Range("W15").Value = 1
$ws.Range("C16").Value = 1
$ws.Range("D16").Value = 1
$ws.Range("E16").Value = 1
$ws.Range("F16").Value = 1
$ws.Range("G16").Value = 1
$ws.Range("H16").Value = 1
$ws.Range("I16").Value = 1
$ws.Range("J16").Value = 1
$ws.Range("K16").Value = 1
$ws.Range("L16").Value = 1
$ws.Range("M16").Value = 1
$ws.Range("N16").Value = 1
$ws.Range("O16").Value = 1
$ws.Range("P16").Value = 1
$ws.Range("Q16").Value = 1
$ws.Range("R16").Value = 1
$ws.Range("S16").Value = 1
$ws.Range("T16").Value = 1
$ws.Range("U16").Value = 1
$ws.Range("V16").Value = 1
$ws.Range("W16").Value = 1
$ws.Range("C17").Value = 1
$ws.Range("D17").Value = 1
$ws.Range("E17").Value = 1
$ws.Range("F17").Value = 1
$ws.Range("G17").Value = 1
$ws.Range("H17").Value = 1
$ws.Range("I17").Value = 1
$ws.Range("J17").Value = 1
$ws.Range("K17").Value = 1
$ws.Range("L17").Value = 1
$ws.Range("M17").Value = 1
$ws.Range("N17").Value = 1
$ws.Range("O17").Value = 1
$ws.Range("P17").Value = 1
$ws.Range("Q17").Value = 1
$ws.Range("R17").Value = 1
$ws.Range("S17").Value = 1
$ws.Range("T17").Value = 1
$ws.Range("U17").Value = 1
$ws.Range("V17").Value = 1
$ws.Range("W17").Value = 1
$ws.Range("C18").Value = 1
$ws.Range("D18").Value = 1
$ws.Range("E18").Value = 1
$ws.Range("F18").Value = 1
$ws.Range("G18").Value = 1
$ws.Range("H18").Value = 1
$ws.Range("I18").Value = 1
$ws.Range("J18").Value = 1
$ws.Range("K18").Value = 1
$ws.Range("L18").Value = 1
$ws.Range("M18").Value = 1
$ws.Range("N18").Value = 1
$ws.Range("O18").Value = 1
$ws.Range("P18").Value = 1
$ws.Range("Q18").Value = 1
$ws.Range("R18").Value = 1
$ws.Range("S18").Value = 1
$ws.Range("T18").Value = 1
$ws.Range("U18").Value = 1
$ws.Range("V18").Value = 1
$ws.Range("W18").Value = 1
$ws.Range("C19").Value = 1
$ws.Range("D19").Value = 1
$ws.Range("E19").Value = 1
$ws.Range("F19").Value = 1
$ws.Range("G19").Value = 1
$ws.Range("H19").Value = 1
$ws.Range("I19").Value = 1
$ws.Range("J19").Value = 1
$ws.Range("K19").Value = 1
$ws.Range("L19").Value = 1
$ws.Range("M19").Value = 1
$ws.Range("N19").Value = 1
$ws.Range("O19").Value = 1
$ws.Range("P19").Value = 1
$ws.Range("Q19").Value = 1
$ws.Range("R19").Value = 1
$ws.Range("S19").Value = 1
$ws.Range("T19").Value = 1
$ws.Range("U19").Value = 1
$ws.Range("V19").Value = 1
$ws.Range("W19").Value = 1
$ws.Range("C20").Value = 1
$ws.Range("D20").Value = 1
$ws.Range("E20").Value = 1
$ws.Range("F20").Value = 1
$ws.Range("G20").Value = 1
$ws.Range("H20").Value = 1
$ws.Range("I20").Value = 1
$ws.Range("J20").Value = 1
$ws.Range("K20").Value = 1
$ws.Range("L20").Value = 1
$ws.Range("M20").Value = 1
$ws.Range("N20").Value = 1
$ws.Range("O20").Value = 1
$ws.Range("P20").Value = 1
$ws.Range("Q20").Value = 1
$ws.Range("R20").Value = 1
$ws.Range("S20").Value = 1
$ws.Range("T20").Value = 1
$ws.Range("U20").Value = 1
$ws.Range("V20").Value = 1
$ws.Range("W20").Value = 1
$ws.Range("C21").Value = 1
$ws.Range("D21").Value = 1
$ws.Range("E21").Value = 1
$ws.Range("F21").Value = 1
$ws.Range("G21").Value = 1
$ws.Range("H21").Value = 1
$ws.Range("I21").Value = 1
$ws.Range("J21").Value = 1
$ws.Range("K21").Value = 1
$ws.Range("L21").Value = 1
$ws.Range("M21").Value = 1
$ws.Range("N21").Value = 1
$ws.Range("O21").Value = 1
$ws.Range("P21").Value = 1
$ws.Range("Q21").Value = 1
$ws.Range("R21").Value = 1
$ws.Range("S21").Value = 1
$ws.Range("T21").Value = 1
$ws.Range("U21").Value = 1
$ws.Range("V21").Value = 1
$ws.Range("W21").Value = 1
$ws.Range("C22").Value = 1
$ws.Range("D22").Value = 1
$ws.Range("E22").Value = 1
$ws.Range("F22").Value = 1
$ws.Range("G22").Value = 1
$ws.Range("H22").Value = 1
$ws.Range("I22").Value = 1
$ws.Range("J22").Value = 1
$ws.Range("K22").Value = 1
$ws.Range("L22").Value = 1
$ws.Range("M22").Value = 1
$ws.Range("N22").Value = 1
$ws.Range("O22").Value = 1
$ws.Range("P22").Value = 1
$ws.Range("Q22").Value = 1
$ws.Range("R22").Value = 1
$ws.Range("S22").Value = 1
$ws.Range("T22").Value = 1
$ws.Range("U22").Value = 1
$ws.Range("V22").Value = 1
$ws.Range("W22").Value = 1
$ws.Range("C23").Value = 1
$ws.Range("D23").Value = 1
$ws.Range("E23").Value = 1
$ws.Range("F23").Value = 1
$ws.Range("G23").Value = 1
$ws.Range("H23").Value = 1
$ws.Range("I23").Value = 1
$ws.Range("J23").Value = 1
$ws.Range("K23").Value = 1
$ws.Range("L23").Value = 1
$ws.Range("M23").Value = 1
$ws.Range("N23").Value = 1
$ws.Range("O23").Value = 1
$ws.Range("P23").Value = 1
$ws.Range("Q23").Value = 1
$ws.Range("R23").Value = 1
$ws.Range("S23").Value = 1
$ws.Range("T23").Value = 1
$ws.Range("U23").Value = 1
$ws.Range("V23").Value = 1
$ws.Range("W23").Value = 1
$ws.Range("C24").Value = 1
$ws.Range("D24").Value = 1
$ws.Range("E24").Value = 1
$ws.Range("F24").Value = 1
$ws.Range("G24").Value = 1
$ws.Range("H24").Value = 1
$ws.Range("I24").Value = 1
$ws.Range("J24").Value = 1
$ws.Range("K24").Value = 1
$ws.Range("L24").Value = 1
$ws.Range("M24").Value = 1
$ws.Range("N24").Value = 1
$ws.Range("O24").Value = 1
$ws.Range("P24").Value = 1
$ws.Range("Q24").Value = 1
$ws.Range("R24").Value = 1
$ws.Range("S24").Value = 1
$ws.Range("T24").Value = 1
$ws.Range("U24").Value = 1
$ws.Range("V24").Value = 1
$ws.Range("W24").Value = 1
$ws.Range("C25").Value = 1
$ws.Range("D25").Value = 1
$ws.Range("E25").Value = 1
$ws.Range("F25").Value = 1
$ws.Range("G25").Value = 1
$ws.Range("H25").Value = 1
$ws.Range("I25").Value = 1
$ws.Range("J25").Value = 1
$ws.Range("K25").Value = 1
$ws.Range("L25").Value = 1
$ws.Range("M25").Value = 1
$ws.Range("N25").Value = 1
$ws.Range("O25").Value = 1
$ws.Range("P25").Value = 1
$ws.Range("Q25").Value = 1
$ws.Range("R25").Value = 1
$ws.Range("S25").Value = 1
$ws.Range("T25").Value = 1
$ws.Range("U25").Value = 1
$ws.Range("V25").Value = 1
$ws.Range("W25").Value = 1
$ws.Range("C26").Value = 1
$ws.Range("D26").Value = 1
$ws.Range("E26").Value = 1
$ws.Range("F26").Value = 1
$ws.Range("G26").Value = 1
$ws.Range("H26").Value = 1
$ws.Range("I26").Value = 1
$ws.Range("J26").Value = 1
$ws.Range("K26").Value = 1
$ws.Range("L26").Value = 1
$ws.Range("M26").Value = 1
$ws.Range("N26").Value = 1
$ws.Range("O26").Value = 1
$ws.Range("P26").Value = 1
$ws.Range("Q26").Value = 1
$ws.Range("R26").Value = 1
$ws.Range("S26").Value = 1
$ws.Range("T26").Value = 1
$ws.Range("U26").Value = 1
$ws.Range("V26").Value = 1
$ws.Range("W26").Value = 1
$ws.Range("C27").Value = 1
$ws.Range("D27").Value = 1
$ws.Range("E27").Value = 1
$ws.Range("F27").Value = 1
$ws.Range("G27").Value = 1
$ws.Range("H27").Value = 1
$ws.Range("I27").Value = 1
$ws.Range("J27").Value = 1
$ws.Range("K27").Value = 1
$ws.Range("L27").Value = 1
$ws.Range("M27").Value = 1
$ws.Range("N27").Value = 1
$ws.Range("O27").Value = 1
$ws.Range("P27").Value = 1
$ws.Range("Q27").Value = 1
$ws.Range("R27").Value = 1
$ws.Range("S27").Value = 1
$ws.Range("T27").Value = 1
$ws.Range("U27").Value = 1
$ws.Range("V27").Value = 1
$ws.Range("W27").Value = 1
$ws.Range("C28").Value = 1
$ws.Range("D28").Value = 1
$ws.Range("E28").Value = 1
$ws.Range("F28").Value = 1
$ws.Range("G28").Value = 1
$ws.Range("H28").Value = 1
$ws.Range("I28").Value = 1
$ws.Range("J28").Value = 1
$ws.Range("K28").Value = 1
$ws.Range("L28").Value = 1
$ws.Range("M28").Value = 1
$ws.Range("N28").Value = 1
$ws.Range("O28").Value = 1
$ws.Range("P28").Value = 1
$ws.Range("Q28").Value = 1
$ws.Range("R28").Value = 1
$ws.Range("S28").Value = 1
$ws.Range("T28").Value = 1
$ws.Range("U28").Value = 1
$ws.Range("V28").Value = 1
$ws.Range("W28").Value = 1
$ws.Range("C29").Value = 1
$ws.Range("D29").Value = 1
$ws.Range("E29").Value = 1
$ws.Range("F29").Value = 1
$ws.Range("G29").Value = 1
$ws.Range("H29").Value = 1
$ws.Range("I29").Value = 1
$ws.Range("J29").Value = 1
$ws.Range("K29").Value = 1
$ws.Range("L29").Value = 1
$ws.Range("M29").Value = 1
$ws.Range("N29").Value = 1
$ws.Range("O29").Value = 1
$ws.Range("P29").Value = 1
$ws.Range("Q29").Value = 1
$ws.Range("R29").Value = 1
$ws.Range("S29").Value = 1
$ws.Range("T29").Value = 1
$ws.Range("U29").Value = 1
$ws.Range("V29").Value = 1
$ws.Range("W29").Value = 1

# --- Styling (bold, thin border, centered/top - matches original style index 1) ---
$headerRange = $ws.Range("B1:W1")
$headerRange.Font.Bold = $true
$headerRange.Borders.LineStyle = 1
$headerRange.HorizontalAlignment = -4108
$headerRange.VerticalAlignment = -4160

$idxRange = $ws.Range("A2:A29")
$idxRange.Font.Bold = $true
$idxRange.Borders.LineStyle = 1
$idxRange.HorizontalAlignment = -4108
$idxRange.VerticalAlignment = -4160

